# Update the ITR sheet task-assignment workbook: revised time estimates (hours -> days)
# and selection/view state, per "Revised Task Assignments and Iteration plan for E1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITR")
$ws.Activate()

# Column H holds the long task-description text for rows 2-9; the revision updates the
# "time to complete" figures at the end of each description from hours to days.
$ws.Range("H5").Value = "Създаване на `"Списък на рисковете`" - подробен списък с рисковете при изпълнението на проекта. Темплейт на документа - `"03-Risk-List-tmp.xls`". Също така, да се създаде и да се разпишат всички точки от документа `"План за управление на рисковете и проблемите`". Темплейт - `"02-Risk-Management-Plan-rup_rskpln.dot`". Време за изпълнение: 3 дни"
$ws.Range("H6").Value = "Създаване на документа `"Речник`". Темплейт - `"rup_gloss.dot`". Подробно описание на всички точки на документа. След създаването му, документът е подготвен и отворен към всички участници в разработката, които имат необходимостта да го ползват/допълват. - 2 дни"
$ws.Range("H9").Value = "Създаване на документ `"График за изпълнението на проекта`". Темплейт - `"07-0-1-Example-Project-Schedule.mpp`".  Да се създаде подробен график по изпълнението, като се разпишат фазите, итерациите и задачите, както и сроковете за изпълнение. Време за изпълнение: 4 дни"
$ws.Range("H7").Value = "Създаване на документа `"Спецификация на допълнителните изисквания`" (SRS). Да се разпишат подробно функционалностите на разработваната система. Да се групират и обособят, като се изведат модулите на системата. Да се попълнят и подробно разпишат всички точки на споменатия документ. Темплейт - `"02-Requirements-Management-Plan-rup_rmpln.dot`". Време за изпълнение: 2 дни"
$ws.Range("H2").Value = "Създаване на `"План за изпълнение на итерация E1 от общия план`". Да се създаде документ подобен на текущия (или ако имаме възможност в Jira), който да дава подробна информация относно изпълнението на следващата итерация от проекта - Е1. Да се разпишат и възложат необходимите задачи, като се определи време за изпълнение и се запише подробно описание на типа на задачата и начина, по който да бъдат представени резултатите. Определено време за изпълнение - 3 дни"
$ws.Range("H8").Value = "Създаване на документа `"Конвенции за писане на код`". Темплейт - `"rup_prggd.dot`". Подробно разписване и попълване на всички точки от документа. Определено време за изпълнение - около 3 дни"
$ws.Range("H4").Value = "Създаване на `"План за управление на качеството`". Темплейт - `"02-Quality-Assurance-Plan-rup_qapln`". Подробно разписване и попълване на всички точки от документа. Определено време за изпълнение - около 4 дни"
$ws.Range("H3").Value = "Създаване на документа `"План за разработка на софтуерния продукт`". Темплейт - `"01-rup_sdpln.dot`". Подробно разписване и попълване на всички точки от документа. Определено време за изпълнение - около 4 дни"

# Update the saved cursor/selection position shown when the workbook is reopened.
$ws.Range("H9").Select()
